$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (D:E) for the two newest quarters; existing D:K shifts to F:M
$ws.Columns("D:E").Insert()

# Copy number formatting from column F (the old column D, now shifted) onto the new D:E columns
$ws.Columns("F").Copy()
$ws.Columns("D:E").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new quarter columns with their reported figures
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43371
$ws.Range("D8").Value = 156200
$ws.Range("E8").Value = 160800
$ws.Range("D9").Value = 91700
$ws.Range("E9").Value = 91200
$ws.Range("D10").Value = 64500
$ws.Range("E10").Value = 69600
$ws.Range("D12").Value = 13300
$ws.Range("E12").Value = 13200
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 3200
$ws.Range("E14").Value = 2300
$ws.Range("D15").Value = 4000
$ws.Range("E15").Value = 3900
$ws.Range("D17").Value = 140500
$ws.Range("E17").Value = 139800
$ws.Range("D18").Value = 15700
$ws.Range("E18").Value = 21000
$ws.Range("D20").Value = -2100
$ws.Range("E20").Value = -2400
$ws.Range("D21").Value = 23200
$ws.Range("E21").Value = 27800
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 13600
$ws.Range("E23").Value = 18600
$ws.Range("D24").Value = 1900
$ws.Range("E24").Value = 3600
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 11600
$ws.Range("E26").Value = 15000
$ws.Range("D27").Value = 11600
$ws.Range("E27").Value = 14600
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 2100
$ws.Range("E32").Value = 2400
$ws.Range("D33").Value = 11600
$ws.Range("E33").Value = 14600
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 11600
$ws.Range("E35").Value = 14600
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43371
$ws.Range("D41").Value = 82000
$ws.Range("E41").Value = 111800
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 85800
$ws.Range("E43").Value = 92900
$ws.Range("D44").Value = 104800
$ws.Range("E44").Value = 98900
$ws.Range("D45").Value = 9200
$ws.Range("E45").Value = 9500
$ws.Range("D46").Value = 281800
$ws.Range("E46").Value = 313200
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 65500
$ws.Range("E48").Value = 66200
$ws.Range("D49").Value = 360600
$ws.Range("E49").Value = 370400
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 11800
$ws.Range("E52").Value = 9600
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 719600
$ws.Range("E54").Value = 759400
$ws.Range("D57").Value = 50700
$ws.Range("E57").Value = 48300
$ws.Range("D58").Value = 4500
$ws.Range("E58").Value = 9100
$ws.Range("D59").Value = 48900
$ws.Range("E59").Value = 51100
$ws.Range("D60").Value = 104200
$ws.Range("E60").Value = 108500
$ws.Range("D61").Value = 202800
$ws.Range("E61").Value = 247300
$ws.Range("D62").Value = 44300
$ws.Range("E62").Value = 44000
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 351300
$ws.Range("E66").Value = 399800
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -79100
$ws.Range("E72").Value = -90700
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 368300
$ws.Range("E76").Value = 359700
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43371
$ws.Range("D81").Value = 11600
$ws.Range("E81").Value = 14600
$ws.Range("D83").Value = 9700
$ws.Range("E83").Value = 9200
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 21900
$ws.Range("E89").Value = 27400
$ws.Range("D91").Value = -400
$ws.Range("E91").Value = 6000
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -3300
$ws.Range("E94").Value = -7600
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -47900
$ws.Range("E100").Value = -15000
$ws.Range("D101").Value = -500
$ws.Range("E101").Value = -200
$ws.Range("D102").Value = -29800
$ws.Range("E102").Value = 4600

# A handful of prior-quarter figures were restated along with this update
$ws.Range("I9").Value = 172400
$ws.Range("I10").Value = -26100
$ws.Range("I17").Value = 133900
$ws.Range("I18").Value = 12400
$ws.Range("I20").Value = -2900
$ws.Range("I32").Value = 2900
